# Generate Report for Handoff
# The b449e6a1-ccaf-46f7-9753-88a10d2650d6.md file has finished its local
# handback cycle and is now ready to be handed off again. Update its
# status everywhere it's reported, and record the new handoff datetimes
# for each target locale.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"

# --- Overview sheet: update status for the b449... file row (row 3) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusReady
$overview.Range("C3").Value = $statusReady

# --- zh-cn sheet: update status + latest handoff datetime (row 3) ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $statusReady
$zhcn.Range("D3").Value = "2016-03-04 03:27:52"

# --- de-de sheet: update status + latest handoff datetime (row 3) ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $statusReady
$dede.Range("D3").Value = "2016-03-04 03:28:08"
